$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text without altering its cell style,
# even when the text looks numeric (e.g. "316.01", "28.119.40", "1.909.95").
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '28.119.40'
Set-TextValue 'E2' '  +2.28%  '
Set-TextValue 'D3' '1.910.09'
Set-TextValue 'E3' '  +1.95%  '
Set-TextValue 'E4' '  -1.15%  '
Set-TextValue 'D5' '316.01'
Set-TextValue 'E5' '  +0.84%  '
Set-TextValue 'E6' '  -1.13%  '
Set-TextValue 'D7' '0.4831'
Set-TextValue 'E7' '  +0.97%  '
Set-TextValue 'D8' '0.3820'
Set-TextValue 'E8' '  +1.24%  '
Set-TextValue 'D9' '0.07355'
Set-TextValue 'E9' '  -0.30%  '
Set-TextValue 'D10' '0.9366'
Set-TextValue 'E10' '  -0.10%  '
Set-TextValue 'D11' '20.80'
Set-TextValue 'E11' '  +0.36%  '
Set-TextValue 'D12' '0.07793'
Set-TextValue 'E12' '  -0.64%  '
Set-TextValue 'D13' '1.917.30'
Set-TextValue 'E13' '  +2.22%  '
Set-TextValue 'D14' '5.517'
Set-TextValue 'E14' '  +1.29%  '
Set-TextValue 'D15' '6.631'
Set-TextValue 'E15' '  +0.60%  '
Set-TextValue 'D16' '91.44'
Set-TextValue 'E16' '  +0.56%  '
Set-TextValue 'D17' '1.005'
Set-TextValue 'E17' '  -1.10%  '
Set-TextValue 'D18' '0.000008832'
Set-TextValue 'E18' '  -1.08%  '
Set-TextValue 'E19' '  -1.04%  '
Set-TextValue 'D20' '28.140.00'
Set-TextValue 'E20' '  +2.21%  '
Set-TextValue 'D21' '14.85'
Set-TextValue 'E21' '  -0.42%  '
Set-TextValue 'D22' '5.145'
Set-TextValue 'E22' '  +0.19%  '
Set-TextValue 'D23' '2.151.44'
Set-TextValue 'E23' '  +1.40%  '
Set-TextValue 'D24' '10.90'
Set-TextValue 'E24' '  +1.59%  '
Set-TextValue 'D25' '156.74'
Set-TextValue 'E25' '  +1.81%  '
Set-TextValue 'E26' '  -1.76%  '
Set-TextValue 'D27' '18.57'
Set-TextValue 'E27' '  +0.06%  '
Set-TextValue 'D28' '2.111'
Set-TextValue 'E28' '  +4.58%  '
Set-TextValue 'D29' '116.37'
Set-TextValue 'E29' '  +0.35%  '
Set-TextValue 'D30' '4.951'
Set-TextValue 'E30' '  -0.93%  '
Set-TextValue 'D31' '0.08922'
Set-TextValue 'E31' '  -0.15%  '
Set-TextValue 'D32' '3.346'
Set-TextValue 'E32' '  +0.01%  '
Set-TextValue 'D33' '1.254'
Set-TextValue 'E33' '  +2.63%  '
Set-TextValue 'D34' '0.7694'
Set-TextValue 'E34' '  +2.18%  '
Set-TextValue 'D35' '4.683'
Set-TextValue 'E35' '  +1.55%  '
Set-TextValue 'D36' '2.611'
Set-TextValue 'E36' '  -2.95%  '
Set-TextValue 'D37' '0.02056'
Set-TextValue 'E37' '  -0.30%  '
Set-TextValue 'E38' '  -1.42%  '
Set-TextValue 'D39' '0.05310'
Set-TextValue 'E39' '  +0.09%  '
Set-TextValue 'D40' '0.5501'
Set-TextValue 'E40' '  +2.73%  '
Set-TextValue 'D41' '2.975'
Set-TextValue 'E41' '  -1.01%  '
Set-TextValue 'D42' '7.027'
Set-TextValue 'D43' '0.1524'
Set-TextValue 'E43' '  -0.26%  '
Set-TextValue 'D44' '8.456'
Set-TextValue 'E44' '  +0.27%  '
Set-TextValue 'D45' '10.66'
Set-TextValue 'E45' '  +0.61%  '
Set-TextValue 'D46' '0.4836'
Set-TextValue 'E46' '  +0.01%  '
Set-TextValue 'D47' '107.40'
Set-TextValue 'E47' '  +4.22%  '
Set-TextValue 'E48' '  -1.19%  '
Set-TextValue 'D49' '1.659'
Set-TextValue 'E49' '  -0.23%  '
Set-TextValue 'D50' '68.38'
Set-TextValue 'E50' '  +1.65%  '
Set-TextValue 'D51' '0.06107'
Set-TextValue 'E51' '  +0.20%  '
